# The "Förändrad" column (C) records the date a permit record was last
# changed. The source data was refreshed, moving that date forward one
# day (2023-10-05, serial 45204 -> 2023-10-06, serial 45205) for every
# data row currently stamped with the old date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count - 1 + $usedRange.Row

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45205
    }
}
